$p = $ppt.ActivePresentation
$p.Slides.Item(15).Delete()
